# fix: Scenarios sheet Named Range 100% conversion + scenario value differentiation
#
# 1. Add 12 new workbook-level defined names (Scen_Method{1..4}_{Best,Base,Worst})
#    pointing at the Method-SAM rows (17-20) on the Scenarios sheet.
# 2. Differentiate the Best/Worst formulas for the four Method SAM rows so they
#    are no longer identical to the Base case (+15% / -15%).
# 3. Re-point the Average SAM row (row 22) at the new named ranges instead of
#    the previously hardcoded (and wrong) B18:B21 / C18:C21 / D18:D21 ranges.
# 4. Bump the "Generated" timestamp on the Validation_Log sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Scenarios")

# --- 1. New named ranges -------------------------------------------------
$wb.Names.Add('Scen_Method1_Best',  "='Scenarios'!`$B`$17")
$wb.Names.Add('Scen_Method1_Base',  "='Scenarios'!`$C`$17")
$wb.Names.Add('Scen_Method1_Worst', "='Scenarios'!`$D`$17")
$wb.Names.Add('Scen_Method2_Best',  "='Scenarios'!`$B`$18")
$wb.Names.Add('Scen_Method2_Base',  "='Scenarios'!`$C`$18")
$wb.Names.Add('Scen_Method2_Worst', "='Scenarios'!`$D`$18")
$wb.Names.Add('Scen_Method3_Best',  "='Scenarios'!`$B`$19")
$wb.Names.Add('Scen_Method3_Base',  "='Scenarios'!`$C`$19")
$wb.Names.Add('Scen_Method3_Worst', "='Scenarios'!`$D`$19")
$wb.Names.Add('Scen_Method4_Best',  "='Scenarios'!`$B`$20")
$wb.Names.Add('Scen_Method4_Base',  "='Scenarios'!`$C`$20")
$wb.Names.Add('Scen_Method4_Worst', "='Scenarios'!`$D`$20")

# --- 2. Differentiate Best/Worst formulas for each Method SAM row -------
$ws.Range("B17").Formula = "=SAM*1.15"
$ws.Range("D17").Formula = "=SAM*0.85"

$ws.Range("B18").Formula = "=SAM_Method2*1.15"
$ws.Range("D18").Formula = "=SAM_Method2*0.85"

$ws.Range("B19").Formula = "=SAM_Method3*1.15"
$ws.Range("D19").Formula = "=SAM_Method3*0.85"

$ws.Range("B20").Formula = "=SAM_Method4*1.15"
$ws.Range("D20").Formula = "=SAM_Method4*0.85"

# --- 3. Average SAM now references the new named ranges -----------------
$ws.Range("B22").Formula = "=AVERAGE(Scen_Method1_Best,Scen_Method2_Best,Scen_Method3_Best,Scen_Method4_Best)"
$ws.Range("C22").Formula = "=AVERAGE(Scen_Method1_Base,Scen_Method2_Base,Scen_Method3_Base,Scen_Method4_Base)"
$ws.Range("D22").Formula = "=AVERAGE(Scen_Method1_Worst,Scen_Method2_Worst,Scen_Method3_Worst,Scen_Method4_Worst)"

# --- 4. Bump the "Generated" timestamp on Validation_Log -----------------
$vl = $wb.Worksheets.Item("Validation_Log")
$vl.Range("A2").Value = "Generated: 2025-11-04 19:31"
